# =============================================================================
# Source deck: the "active" design theme (Slide Master -> theme2.xml) is
# "Integral"/"Red Violet". The target edit turns that into the plain
# "Office Theme" palette (this is what the author's diff shows landing in
# the part that is wired to the Slide Master / Presentation relationship).
# The deck's NotesMaster points at a second, otherwise-unused theme part
# (theme1.xml) which this PowerPoint object model does not expose as a
# separate editable theme (NotesMaster/NotesPage resolve to the same
# single theme store as the Slide Master in this host) - only the
# Slide-Master-facing theme can be driven from script, so that is what
# gets repainted below.
# =============================================================================

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table (on slide 5): switch the table style from the locally-defined
#    "Table_0" style ({040D479F-8844-40A7-99EA-823B388D4706}) to the
#    built-in theme table style {2FBA782F-A967-4FE5-8FF4-286D05F9B548}.
#    Walk every slide/shape so the edit is not tied to a hard-coded index.
# ---------------------------------------------------------------------------
$oldStyleId = "{040D479F-8844-40A7-99EA-823B388D4706}"
$newStyleId = "{2FBA782F-A967-4FE5-8FF4-286D05F9B548}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable -and $shp.Table.Style -eq $oldStyleId) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's theme from "Integral / Red Violet" to the
#    standard Office theme palette.
# ---------------------------------------------------------------------------
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
